$wb = $excel.ActiveWorkbook

# Delete the placeholder header row (row 1) from the "Steps 3 & 4" sheet.
# This promotes the real header row (old row 2) to row 1, and all data
# rows shift up by one.
$ws3 = $wb.Worksheets.Item("Steps 3 & 4")
$ws3.Rows.Item(1).Delete()

# Make "Steps 3 & 4" the active/selected sheet (it becomes tab index 3,
# i.e. activeTab = 2 zero-based).
$ws3.Activate()
